$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text columns (B,D,E,F,G,H,I,J) keep their text type
# by pre-formatting the full data range as Text before assigning string values.
# (Must be one contiguous range - multi-area comma-joined ranges only apply
# NumberFormat to the first area in this COM host.)
$ws.Range("B1123:J1154").NumberFormat = "@"

$ws.Range("C1123").Value = "Arizona Cardinals"
$ws.Range("D1123").Value = "790"
$ws.Range("E1123").Value = "558"
$ws.Range("F1123").Value = "232"
$ws.Range("G1123").Value = "29.0"
$ws.Range("H1123").Value = "13"
$ws.Range("I1123").Value = "6"

$ws.Range("C1124").Value = "Kansas City Chiefs"
$ws.Range("D1124").Value = "733"
$ws.Range("E1124").Value = "501"
$ws.Range("F1124").Value = "232"
$ws.Range("G1124").Value = "32.0"
$ws.Range("I1124").Value = "6"
$ws.Range("J1124").Value = "2"

$ws.Range("C1125").Value = "Pittsburgh Steelers"
$ws.Range("D1125").Value = "727"
$ws.Range("E1125").Value = "528"
$ws.Range("F1125").Value = "199"
$ws.Range("G1125").Value = "34.0"
$ws.Range("H1125").Value = "13"
$ws.Range("I1125").Value = "12"
$ws.Range("J1125").Value = "2"

$ws.Range("C1126").Value = "Washington Redskins"
$ws.Range("D1126").Value = "721"
$ws.Range("E1126").Value = "496"
$ws.Range("F1126").Value = "225"
$ws.Range("G1126").Value = "23.0"
$ws.Range("H1126").Value = "7"
$ws.Range("I1126").Value = "4"

$ws.Range("C1127").Value = "Los Angeles Chargers"
$ws.Range("D1127").Value = "720"
$ws.Range("E1127").Value = "494"
$ws.Range("F1127").Value = "226"
$ws.Range("G1127").Value = "25.0"
$ws.Range("H1127").Value = "6"
$ws.Range("I1127").Value = "3"
$ws.Range("J1127").Value = "0"

$ws.Range("C1128").Value = "Los Angeles Rams"
$ws.Range("D1128").Value = "719"
$ws.Range("E1128").Value = "464"
$ws.Range("F1128").Value = "255"
$ws.Range("G1128").Value = "29.0"
$ws.Range("H1128").Value = "13"
$ws.Range("I1128").Value = "8"
$ws.Range("J1128").Value = "2"

$ws.Range("C1129").Value = "New York Giants"
$ws.Range("D1129").Value = "705"
$ws.Range("E1129").Value = "457"
$ws.Range("F1129").Value = "248"
$ws.Range("G1129").Value = "24.0"
$ws.Range("H1129").Value = "5"
$ws.Range("I1129").Value = "4"

$ws.Range("C1130").Value = "Houston Texans"
$ws.Range("D1130").Value = "704"
$ws.Range("E1130").Value = "503"
$ws.Range("F1130").Value = "201"
$ws.Range("H1130").Value = "13"
$ws.Range("I1130").Value = "9"

$ws.Range("C1131").Value = "Atlanta Falcons"
$ws.Range("D1131").Value = "699"
$ws.Range("E1131").Value = "448"
$ws.Range("F1131").Value = "251"
$ws.Range("G1131").Value = "18.0"
$ws.Range("H1131").Value = "8"
$ws.Range("I1131").Value = "2"
$ws.Range("J1131").Value = "0"

$ws.Range("C1132").Value = "Cincinnati Bengals"
$ws.Range("D1132").Value = "698"
$ws.Range("E1132").Value = "442"
$ws.Range("F1132").Value = "256"
$ws.Range("G1132").Value = "13.0"
$ws.Range("H1132").Value = "6"
$ws.Range("I1132").Value = "5"

$ws.Range("C1133").Value = "Tennessee Titans"
$ws.Range("D1133").Value = "697"
$ws.Range("E1133").Value = "448"
$ws.Range("F1133").Value = "249"
$ws.Range("G1133").Value = "27.0"
$ws.Range("H1133").Value = "7"
$ws.Range("I1133").Value = "6"
$ws.Range("J1133").Value = "1"

$ws.Range("C1134").Value = "Indianapolis Colts"
$ws.Range("D1134").Value = "689"
$ws.Range("E1134").Value = "513"
$ws.Range("F1134").Value = "176"
$ws.Range("G1134").Value = "26.0"
$ws.Range("H1134").Value = "10"
$ws.Range("I1134").Value = "5"

$ws.Range("C1135").Value = "Miami Dolphins"
$ws.Range("D1135").Value = "686"
$ws.Range("E1135").Value = "444"
$ws.Range("F1135").Value = "242"
$ws.Range("G1135").Value = "13.0"
$ws.Range("I1135").Value = "1"

$ws.Range("C1136").Value = "New York Jets"
$ws.Range("D1136").Value = "679"
$ws.Range("E1136").Value = "482"
$ws.Range("F1136").Value = "197"
$ws.Range("G1136").Value = "25.0"
$ws.Range("I1136").Value = "7"
$ws.Range("J1136").Value = "2"

$ws.Range("C1137").Value = "Minnesota Vikings"
$ws.Range("D1137").Value = "672"
$ws.Range("E1137").Value = "494"
$ws.Range("F1137").Value = "178"
$ws.Range("G1137").Value = "31.0"
$ws.Range("H1137").Value = "8"
$ws.Range("I1137").Value = "6"
$ws.Range("J1137").Value = "0"

$ws.Range("D1138").Value = "666"
$ws.Range("E1138").Value = "507"
$ws.Range("F1138").Value = "159"
$ws.Range("G1138").Value = "20.0"
$ws.Range("H1138").Value = "14"
$ws.Range("I1138").Value = "10"

$ws.Range("B1139").Value = "17"
$ws.Range("C1139").Value = "Carolina Panthers"
$ws.Range("D1139").Value = "657"
$ws.Range("E1139").Value = "461"
$ws.Range("F1139").Value = "196"
$ws.Range("G1139").Value = "39.0"
$ws.Range("I1139").Value = "7"
$ws.Range("J1139").Value = "1"

$ws.Range("B1140").Value = "18"
$ws.Range("C1140").Value = "Seattle Seahawks"
$ws.Range("D1140").Value = "650"
$ws.Range("E1140").Value = "412"
$ws.Range("F1140").Value = "238"
$ws.Range("G1140").Value = "20.0"
$ws.Range("H1140").Value = "11"
$ws.Range("I1140").Value = "11"
$ws.Range("J1140").Value = "1"

$ws.Range("C1141").Value = "Buffalo Bills"
$ws.Range("D1141").Value = "647"
$ws.Range("E1141").Value = "449"
$ws.Range("F1141").Value = "198"
$ws.Range("G1141").Value = "29.0"
$ws.Range("H1141").Value = "9"
$ws.Range("I1141").Value = "6"

$ws.Range("C1142").Value = "Green Bay Packers"
$ws.Range("D1142").Value = "646"
$ws.Range("E1142").Value = "437"
$ws.Range("F1142").Value = "209"
$ws.Range("G1142").Value = "25.0"
$ws.Range("H1142").Value = "9"
$ws.Range("I1142").Value = "6"
$ws.Range("J1142").Value = "0"

$ws.Range("C1143").Value = "Dallas Cowboys"
$ws.Range("D1143").Value = "640"
$ws.Range("E1143").Value = "428"
$ws.Range("F1143").Value = "212"
$ws.Range("G1143").Value = "26.0"
$ws.Range("H1143").Value = "10"
$ws.Range("I1143").Value = "8"

$ws.Range("C1144").Value = "Chicago Bears"
$ws.Range("D1144").Value = "639"
$ws.Range("E1144").Value = "498"
$ws.Range("F1144").Value = "141"
$ws.Range("G1144").Value = "25.0"
$ws.Range("H1144").Value = "11"
$ws.Range("I1144").Value = "7"

$ws.Range("C1145").Value = "Denver Broncos"
$ws.Range("D1145").Value = "628"
$ws.Range("E1145").Value = "441"
$ws.Range("F1145").Value = "187"
$ws.Range("G1145").Value = "24.0"
$ws.Range("H1145").Value = "7"
$ws.Range("I1145").Value = "3"
$ws.Range("J1145").Value = "0"

$ws.Range("C1146").Value = "Tampa Bay Buccaneers"
$ws.Range("D1146").Value = "625"
$ws.Range("E1146").Value = "456"
$ws.Range("F1146").Value = "169"
$ws.Range("G1146").Value = "22.0"
$ws.Range("H1146").Value = "11"
$ws.Range("I1146").Value = "9"
$ws.Range("J1146").Value = "1"

$ws.Range("C1147").Value = "Oakland Raiders"
$ws.Range("D1147").Value = "623"
$ws.Range("E1147").Value = "466"
$ws.Range("F1147").Value = "157"
$ws.Range("G1147").Value = "25.0"
$ws.Range("H1147").Value = "6"
$ws.Range("I1147").Value = "4"
$ws.Range("J1147").Value = "0"

$ws.Range("C1148").Value = "Cleveland Browns"
$ws.Range("D1148").Value = "621"
$ws.Range("E1148").Value = "468"
$ws.Range("F1148").Value = "153"
$ws.Range("G1148").Value = "30.0"
$ws.Range("I1148").Value = "5"

$ws.Range("C1149").Value = "Jacksonville Jaguars"
$ws.Range("D1149").Value = "612"
$ws.Range("E1149").Value = "433"
$ws.Range("F1149").Value = "179"
$ws.Range("G1149").Value = "33.0"
$ws.Range("H1149").Value = "8"

$ws.Range("C1150").Value = "San Francisco 49ers"
$ws.Range("D1150").Value = "604"
$ws.Range("E1150").Value = "428"
$ws.Range("F1150").Value = "176"
$ws.Range("G1150").Value = "39.0"
$ws.Range("H1150").Value = "16"
$ws.Range("I1150").Value = "11"
$ws.Range("J1150").Value = "2"

$ws.Range("C1151").Value = "Philadelphia Eagles"
$ws.Range("D1151").Value = "587"
$ws.Range("E1151").Value = "416"
$ws.Range("F1151").Value = "171"
$ws.Range("G1151").Value = "25.0"
$ws.Range("H1151").Value = "9"
$ws.Range("I1151").Value = "5"
$ws.Range("J1151").Value = "1"

$ws.Range("C1152").Value = "New England Patriots"
$ws.Range("D1152").Value = "569"
$ws.Range("E1152").Value = "417"
$ws.Range("F1152").Value = "152"
$ws.Range("G1152").Value = "37.0"
$ws.Range("H1152").Value = "11"
$ws.Range("I1152").Value = "9"
$ws.Range("J1152").Value = "2"

$ws.Range("C1153").Value = "Baltimore Ravens"
$ws.Range("D1153").Value = "558"
$ws.Range("E1153").Value = "408"
$ws.Range("F1153").Value = "150"
$ws.Range("G1153").Value = "23.0"
$ws.Range("H1153").Value = "8"
$ws.Range("I1153").Value = "7"
$ws.Range("J1153").Value = "3"

$ws.Range("C1154").Value = "New Orleans Saints"
$ws.Range("D1154").Value = "547"
$ws.Range("E1154").Value = "404"
$ws.Range("F1154").Value = "143"
$ws.Range("G1154").Value = "27.0"
$ws.Range("H1154").Value = "7"
$ws.Range("I1154").Value = "6"
$ws.Range("J1154").Value = "1"
